# Table 2.1.C update: add November 2016 data and refresh dependent totals.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the title / subtitle text that references the reporting month.
$ws.Range("A2").Value = "by Sector, 2006-November 2016 (Thousand Tons)"

# 2) Insert a new row before row 53 (after October 2016, row 52) for the
#    November 2016 entry. This shifts all subsequent rows down by one,
#    which Excel also takes care of for merged cells / dimension.
$ws.Rows("53:53").Insert()

# Copy the formatting of the October row (52) onto the newly inserted row
# so the new row reuses the same cell styles as the other month rows.
$ws.Range("A52:F52").Copy()
$ws.Range("A53:F53").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the November 2016 values.
$ws.Range("A53").Value = "November"
$ws.Range("B53").Value = 49162
$ws.Range("C53").Value = 35358
$ws.Range("D53").Value = 12744
$ws.Range("E53").Value = 61
$ws.Range("F53").Value = 998

# 3) Refresh the "Year to Date" annual total rows (now rows 55-57) to
#    include the new November figures.
$ws.Range("B55").Value = 802369
$ws.Range("C55").Value = 575150
$ws.Range("D55").Value = 208813
$ws.Range("E55").Value = 973
$ws.Range("F55").Value = 17432

$ws.Range("B56").Value = 704639
$ws.Range("C56").Value = 502565
$ws.Range("D56").Value = 185714
$ws.Range("E56").Value = 726
$ws.Range("F56").Value = 15634

$ws.Range("B57").Value = 626108
$ws.Range("C57").Value = 451884
$ws.Range("D57").Value = 160458
$ws.Range("E57").Value = 621
$ws.Range("F57").Value = 13145

# 4) Refresh the "Rolling 12 Months Ending in ..." header text and its
#    totals (now rows 58-60).
$ws.Range("A58").Value = "Rolling 12 Months Ending in November"

$ws.Range("B59").Value = 774011
$ws.Range("C59").Value = 552627
$ws.Range("D59").Value = 203291
$ws.Range("E59").Value = 816
$ws.Range("F59").Value = 17278

$ws.Range("B60").Value = 677695
$ws.Range("C60").Value = 489857
$ws.Range("D60").Value = 172650
$ws.Range("E60").Value = 693
$ws.Range("F60").Value = 14496

$wb.Save()
